# BOM.xlsx: "changes CAN ports to 3 pin, Thermocouple ports to 2 pin"
#
# Row 13 = CN2,CN7  (CAN ports)           -> switch from 2-pin to 3-pin header
# Row 16 = CN5,CN6  (Thermocouple ports)  -> switch from 3-pin to 2-pin header
#
# Columns D (Mfg Part #), E (Description/Value), F (Package/Footprint) and
# H (DigiKey URL) simply swap between the two rows, since the part that used
# to be in row 16 (3-pin / 705510002) now belongs in row 13, and the part
# that used to be in row 13 (2-pin / 705510001) now belongs in row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- capture the "before" values ----------------------------------------
$d13 = $ws.Range("D13").Value()
$e13 = $ws.Range("E13").Value()
$f13 = $ws.Range("F13").Value()
$h13 = $ws.Range("H13").Value()

$d16 = $ws.Range("D16").Value()
$e16 = $ws.Range("E16").Value()
$f16 = $ws.Range("F16").Value()
$h16 = $ws.Range("H16").Value()

# --- row 13 (CAN ports CN2,CN7): now uses the 3-pin (705510002) header -----
$ws.Range("D13").Value = $d16
$ws.Range("E13").Value = $e16
$ws.Range("F13").Value = $f16
$ws.Range("H13").Value = $h16

# --- row 16 (Thermocouple ports CN5,CN6): now uses the 2-pin (705510001)
#     header. D16 must hold a genuine number (705510001), so the number
#     format is relaxed before the assignment and restored afterwards to
#     avoid the text-format cell coercing it into a text string.
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Value = $d13
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").Value = $e13
$ws.Range("F16").Value = $f13
$ws.Range("H16").Value = $h13

# --- row 4: D4/E4 no longer carry any formatting (cells cleared) ----------
$ws.Range("D4:E4").Clear()

# --- row 56: crimp-socket quantity bumped from 16 to 18 -------------------
$ws.Range("B56").Value = 18

# --- restore the on-screen selection that was active when the author saved
$ws.Range("C57").Select()
